$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45792, "Desayuno", "Astringente", 3, 25773),
    @(45792, "Desayuno", "Blanda", 19, 150176),
    @(45792, "Desayuno", "Coronaria", 15, 128715),
    @(45792, "Desayuno", "Hepatica", 1, 10189),
    @(45792, "Desayuno", "Hipercalorica", 1, 9993),
    @(45792, "Desayuno", "Hiperproteica", 2, 19986),
    @(45792, "Desayuno", "Hipo Grasa", 6, 0),
    @(45792, "Desayuno", "Hipoglucida", 12, 104184),
    @(45792, "Desayuno", "Hiposodica", 28, 240548),
    @(45792, "Desayuno", "Liquida Clara", 4, 41240),
    @(45792, "Desayuno", "Liquida Total", 7, 75957),
    @(45792, "Desayuno", "Liquida Total 140 Cc", 1, 10725),
    @(45792, "Desayuno", "Liquida Total Miel 140 Cc", 2, 21702),
    @(45792, "Desayuno", "Liquida Total Nectar", 9, 97659),
    @(45792, "Desayuno", "Liquida Total Nectar 140 Cc", 2, 21702),
    @(45792, "Desayuno", "Liquida total Miel", 3, 32553),
    @(45792, "Desayuno", "Normal", 57, 450528),
    @(45792, "Desayuno", "Renal Dialisis", 4, 34324),
    @(45792, "Desayuno", "Renal PRE Dialisis", 5, 42905),
    @(45792, "Desayuno", "Semiblanda", 28, 221312),
    @(45792, "Desayuno", "Semiblanda Pequena", 4, 29148),
    @(45792, "Desayuno", "Todo Pure", 2, 23216),
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
